# Fruta / hortaliza, semanal
# Insert 4 new weekly records for "Feria Lagunitas de Puerto Montt - Ciruela"
# right before the existing row 167, pushing the remaining history down
# (old row 167 -> new row 171, ..., old row 234 -> new row 238), then
# populate the 4 newly-opened rows (167-170) with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above row 167 (shifts 167:234 down to 171:238).
$ws.Rows("167:170").Insert()

# Helper to fill one data row with the 20 standard columns (positional args
# only -- named parameter binding is not reliable in this interpreter).
function Set-CiruelaRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad) {
    $ws.Cells.Item($Row, 1).Value = 4
    $ws.Cells.Item($Row, 2).Value = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($Row, 3).Value = "Los Lagos"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 10
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value = 100103
    $ws.Cells.Item($Row, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($Row, 9).Value = 100103002
    $ws.Cells.Item($Row, 10).Value = "Ciruela"
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

# Row 167: Angeleno / Primera
Set-CiruelaRow 167 44726 "Angeleno" "Primera" 400 14000 15000 14500 "`$/caja 15 kilos granel" "Región de O'Higgins" 967 15

# Row 168: Angeleno / Segunda
Set-CiruelaRow 168 44726 "Angeleno" "Segunda" 200 13000 13000 13000 "`$/caja 15 kilos granel" "Región de O'Higgins" 867 15

# Row 169: Pink Delight / Primera
Set-CiruelaRow 169 44726 "Pink Delight" "Primera" 300 14000 15000 14500 "`$/caja 15 kilos granel" "Región de O'Higgins" 967 15

# Row 170: Pink Delight / Segunda
Set-CiruelaRow 170 44726 "Pink Delight" "Segunda" 150 13000 13000 13000 "`$/caja 15 kilos granel" "Región de O'Higgins" 867 15
